# Scheduled-runner style market-data refresh: updates cached price/profit
# values (currentAveragePrice*, LevePrice*, LeveProfit*) across the ALC,
# ARM, BSM, CRP, GSM, LTW and WVR leve sheets. Values are plain numbers
# (no formulas back these columns), so each changed cell is written
# directly. A few rows had their HQ price columns drop to 0, which also
# makes the corresponding LeveProfitHQ cell blank (cleared) rather than 0.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 850.06384
$ws.Range("I129").Value = 357.33334
$ws.Range("K129").Value = 1072.00002
$ws.Range("M129").Value = 3927.99998
$ws.Range("H132").Value = 5561265
$ws.Range("J132").Value = 11014.267
$ws.Range("L132").Value = 33042.801
$ws.Range("N132").Value = -38102.801
$ws.Range("H133").Value = 32395.8
$ws.Range("J133").Value = 32395.8
$ws.Range("L133").Value = 32395.8
$ws.Range("N133").Value = -42515.8
$ws.Range("H137").Value = 1064.4861
$ws.Range("I137").Value = 870.8205
$ws.Range("J137").Value = 1293.3636
$ws.Range("K137").Value = 2612.4615
$ws.Range("L137").Value = 3880.0908
$ws.Range("M137").Value = -62.46150000000034
$ws.Range("N137").Value = -8980.0908
$ws.Range("H138").Value = 1184.0928
$ws.Range("I138").Value = 505.2143
$ws.Range("J138").Value = 1702.509
$ws.Range("K138").Value = 1515.6429
$ws.Range("L138").Value = 5107.527
$ws.Range("M138").Value = 3624.3571
$ws.Range("N138").Value = -15387.527
$ws.Range("H139").Value = 34480
$ws.Range("J139").Value = 34480
$ws.Range("L139").Value = 34480
$ws.Range("N139").Value = -44760
$ws.Range("H141").Value = 605.4386
$ws.Range("I141").Value = 509.8113
$ws.Range("J141").Value = 1872.5
$ws.Range("K141").Value = 1529.4339
$ws.Range("L141").Value = 5617.5
$ws.Range("M141").Value = 3650.5661
$ws.Range("N141").Value = -15977.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3315.3425
$ws.Range("I32").Value = 3308.7878
$ws.Range("J32").Value = 3377.1428
$ws.Range("K32").Value = 3308.7878
$ws.Range("L32").Value = 3377.1428
$ws.Range("M32").Value = -3021.7878
$ws.Range("N32").Value = -3951.1428
$ws.Range("H61").Value = 991.4722
$ws.Range("I61").Value = 845.78125
$ws.Range("J61").Value = 2157
$ws.Range("K61").Value = 845.78125
$ws.Range("L61").Value = 2157
$ws.Range("M61").Value = -633.78125
$ws.Range("N61").Value = -2581
$ws.Range("H110").Value = 2464.625
$ws.Range("I110").Value = 1968
$ws.Range("J110").Value = 2762.6
$ws.Range("K110").Value = 1968
$ws.Range("L110").Value = 2762.6
$ws.Range("M110").Value = 77
$ws.Range("N110").Value = -6852.6
$ws.Range("H132").Value = 2341.122
$ws.Range("I132").Value = 2262.8572
$ws.Range("J132").Value = 2509.6924
$ws.Range("K132").Value = 6788.571599999999
$ws.Range("L132").Value = 7529.0772
$ws.Range("M132").Value = -4258.571599999999
$ws.Range("N132").Value = -12589.0772
$ws.Range("H136").Value = 991.4722
$ws.Range("I136").Value = 845.78125
$ws.Range("J136").Value = 2157
$ws.Range("K136").Value = 2537.34375
$ws.Range("L136").Value = 6471
$ws.Range("M136").Value = 12.65625
$ws.Range("N136").Value = -11571
$ws.Range("H139").Value = 28507.5
$ws.Range("J139").Value = 28507.5
$ws.Range("L139").Value = 28507.5
$ws.Range("N139").Value = -38787.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 90910930
$ws.Range("I105").Value = 111112930
$ws.Range("J105").Value = 1955.5
$ws.Range("K105").Value = 111112930
$ws.Range("L105").Value = 1955.5
$ws.Range("M105").Value = -111111183
$ws.Range("N105").Value = -5449.5
$ws.Range("H134").Value = 4797.946
$ws.Range("I134").Value = 1302.1034
$ws.Range("J134").Value = 17470.375
$ws.Range("K134").Value = 3906.3102
$ws.Range("L134").Value = 52411.125
$ws.Range("M134").Value = -1371.3102
$ws.Range("N134").Value = -57481.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 142858640
$ws.Range("I16").Value = 142858640
$ws.Range("K16").Value = 142858640
$ws.Range("M16").Value = -142858353
$ws.Range("H31").Value = 1904.4615
$ws.Range("I31").Value = 1920.64
$ws.Range("K31").Value = 1920.64
$ws.Range("M31").Value = -1625.64
$ws.Range("H34").Value = 1904.4615
$ws.Range("I34").Value = 1920.64
$ws.Range("K34").Value = 1920.64
$ws.Range("M34").Value = -1718.64
$ws.Range("H58").Value = 756.25714
$ws.Range("I58").Value = 684.37036
$ws.Range("J58").Value = 998.875
$ws.Range("K58").Value = 684.37036
$ws.Range("L58").Value = 998.875
$ws.Range("M58").Value = -481.37036
$ws.Range("N58").Value = -1404.875
$ws.Range("H113").Value = 142858640
$ws.Range("I113").Value = 142858640
$ws.Range("K113").Value = 142858640
$ws.Range("M113").Value = -142856470
$ws.Range("H132").Value = 4594.6763
$ws.Range("I132").Value = 5033.6206
$ws.Range("J132").Value = 2048.8
$ws.Range("K132").Value = 15100.8618
$ws.Range("L132").Value = 6146.400000000001
$ws.Range("M132").Value = -12570.8618
$ws.Range("N132").Value = -11206.4
$ws.Range("H136").Value = 756.25714
$ws.Range("I136").Value = 684.37036
$ws.Range("J136").Value = 998.875
$ws.Range("K136").Value = 2053.11108
$ws.Range("L136").Value = 2996.625
$ws.Range("M136").Value = 496.8889199999999
$ws.Range("N136").Value = -8096.625
$ws.Range("H141").Value = 29960
$ws.Range("J141").Value = 29960
$ws.Range("L141").Value = 29960
$ws.Range("N141").Value = -40320

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3255.4443
$ws.Range("I80").Value = 2333.1667
$ws.Range("J80").Value = 5100
$ws.Range("K80").Value = 2333.1667
$ws.Range("L80").Value = 5100
$ws.Range("M80").Value = -1335.1667
$ws.Range("N80").Value = -7096
$ws.Range("H83").Value = 3255.4443
$ws.Range("I83").Value = 2333.1667
$ws.Range("J83").Value = 5100
$ws.Range("K83").Value = 11665.8335
$ws.Range("L83").Value = 25500
$ws.Range("M83").Value = -6673.833500000001
$ws.Range("N83").Value = -35484
$ws.Range("H122").Value = 3951.1667
$ws.Range("I122").Value = 4000
$ws.Range("J122").Value = 3941.4
$ws.Range("K122").Value = 12000
$ws.Range("L122").Value = 11824.2
$ws.Range("M122").Value = -9550
$ws.Range("N122").Value = -16724.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 972.7778
$ws.Range("I16").Value = 969.375
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 969.375
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -799.375
$ws.Range("N16").Value = -1340
$ws.Range("H22").Value = 1860
$ws.Range("I22").Value = 1580
$ws.Range("K22").Value = 1580
$ws.Range("M22").Value = -1285
$ws.Range("H27").Value = 1860
$ws.Range("I27").Value = 1580
$ws.Range("K27").Value = 1580
$ws.Range("M27").Value = -1473
$ws.Range("H40").Value = 4286.4443
$ws.Range("I40").Value = 2367.5715
$ws.Range("J40").Value = 11002.5
$ws.Range("K40").Value = 2367.5715
$ws.Range("L40").Value = 11002.5
$ws.Range("M40").Value = -2231.5715
$ws.Range("N40").Value = -11274.5
$ws.Range("H82").Value = 1915.2084
$ws.Range("I82").Value = 1791
$ws.Range("J82").Value = 2039.4166
$ws.Range("K82").Value = 1791
$ws.Range("L82").Value = 2039.4166
$ws.Range("M82").Value = -1430
$ws.Range("N82").Value = -2761.4166
$ws.Range("H85").Value = 1915.2084
$ws.Range("I85").Value = 1791
$ws.Range("J85").Value = 2039.4166
$ws.Range("K85").Value = 1791
$ws.Range("L85").Value = 2039.4166
$ws.Range("M85").Value = -543
$ws.Range("N85").Value = -4535.4166
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H93").Value = 751.25
$ws.Range("I93").Value = 751.25
$ws.Range("K93").Value = 751.25
$ws.Range("M93").Value = 496.75
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H132").Value = 23804.543
$ws.Range("I132").Value = 1626.174
$ws.Range("J132").Value = 45982.914
$ws.Range("K132").Value = 4878.522
$ws.Range("L132").Value = 137948.742
$ws.Range("M132").Value = -2348.522
$ws.Range("N132").Value = -143008.742

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H96").Value = 4280.143
$ws.Range("I96").Value = 4795.2
$ws.Range("K96").Value = 4795.2
$ws.Range("M96").Value = -3422.2
$ws.Range("H126").Value = 35715492
$ws.Range("I126").Value = 45455476
$ws.Range("J126").Value = 2216.5
$ws.Range("K126").Value = 136366428
$ws.Range("L126").Value = 6649.5
$ws.Range("M126").Value = -136363958
$ws.Range("N126").Value = -11589.5
$ws.Range("H132").Value = 2572.175
$ws.Range("I132").Value = 2514.1143
$ws.Range("K132").Value = 7542.342900000001
$ws.Range("M132").Value = -5012.342900000001
$ws.Range("H136").Value = 668.1923
$ws.Range("I136").Value = 429.41666
$ws.Range("K136").Value = 1288.24998
$ws.Range("M136").Value = 1261.75002
